# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect refreshed data as of commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F
$updates = @{
    2  = 121
    3  = 405
    4  = 11879
    5  = 1177
    6  = 123
    10 = 174
    11 = 100
    13 = 56
    17 = 1439
    18 = 81
    19 = 918
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
